$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update times for rows 3 and 4 (02_data, 03_visualization) from 10.30-12.30 to 10.30-12.00
$ws.Range("B3").Value = "10.30-12.00"
$ws.Range("B4").Value = "10.30-12.00"

# Row 5: Income · Geometries -- add Code/R/RData assignment columns
$ws.Range("E5").Value = "04_income"
$ws.Range("F5").Value = "04_income.R"
$ws.Range("G5").Value = "04_income.RData"

# Row 6: Wealth · Scales -- add Code/R/RData assignment columns (Chart "lines" stays in I6)
$ws.Range("E6").Value = "05_wealth"
$ws.Range("F6").Value = "05_wealth.R"
$ws.Range("G6").Value = "05_wealth.RData"

# Row 7: Mobility · Colors -- add Code assignment column
$ws.Range("E7").Value = "06_mobility"

# Update selection to reflect the new active cell
$ws.Range("E8").Select()
